$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 (subject 27): condition stays "F", fill in Q1..Q15 answers ---
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 3
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 2
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = 2
$ws.Range("M29").Value = 2
$ws.Range("N29").Value = 3
$ws.Range("O29").Value = 3
$ws.Range("P29").Value = 2
$ws.Range("Q29").Value = 4

# --- Row 30 (subject 28): condition changes from "G" to "F", fill in Q1..Q15 answers ---
$ws.Range("B30").Value = "F"
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 3
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 2
$ws.Range("I30").Value = 3
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = 4
$ws.Range("L30").Value = 3
$ws.Range("M30").Value = 4
$ws.Range("N30").Value = 4
$ws.Range("O30").Value = 4
$ws.Range("P30").Value = 3
$ws.Range("Q30").Value = 6

# --- Row 31 (subject 29): condition changes from "F" to "G", fill in Q1..Q15 answers ---
$ws.Range("B31").Value = "G"
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 2
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 3
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 3
$ws.Range("N31").Value = 3
$ws.Range("O31").Value = 4
$ws.Range("P31").Value = 2
$ws.Range("Q31").Value = 3

# --- Row 32 (subject 30): condition changes from "G" to "F" ---
$ws.Range("B32").Value = "F"

# --- Row 33 (new subject 31): add ID + condition "G" ---
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "G"

# --- Update the active selection to match the editor's last position ---
[void]$ws.Range("Q31").Select()
